# Commit message: "adjust all missing value codes to NaN"
#
# The "ColumnHeaders" sheet (the first sheet in the workbook) used the
# literal missing-value code "n/n" (entered with a leading apostrophe,
# i.e. text-quote-prefixed) in cells F29/F30. This edit replaces that
# code with "NaN" (matching the other missing-value-code cells on the
# sheet, e.g. F19:F28 which already read "NaN"), drops the now-orphaned
# "n/n" shared string, and also tidies up a couple of related cosmetic
# details that Excel re-derives when the sheet is touched: the
# dateTimeFormatString note in E6:E8 loses its stray UTC-offset suffix,
# row 23 reverts to the default (auto-fit) row height, and the sheet's
# active selection moves to F27 (where the edit was made).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ColumnHeaders")

# Replace the "n/n" missing-value code with "NaN" in both cells that used
# it. Setting Value2 (rather than Formula/Value) stores a plain literal,
# which also clears the quotePrefix ('n/n was entered as text) formatting
# those two cells previously needed.
$ws.Range("F29").Value2 = "NaN"
$ws.Range("F30").Value2 = "NaN"

# The dateTimeFormatString cells no longer carry a UTC-offset suffix.
$ws.Range("E6").Value2 = "YYYY-MM-DD hh:mm:ss"
$ws.Range("E7").Value2 = "YYYY-MM-DD hh:mm:ss"
$ws.Range("E8").Value2 = "YYYY-MM-DD hh:mm:ss"

# Row 23 goes back to the sheet's default (non-custom) row height.
$ws.Rows.Item(23).AutoFit()

# Reflect where editing finished: cell F27 selected.
$ws.Activate()
$ws.Range("F27").Select()
